# The two "TV4 / S3 / PYR" rows (worksheet rows 6 and 7) are removed.
# This shifts the last row (TV1 / S5 / A14, formerly row 8) up to row 6,
# and the engine automatically drops the now-unused shared strings
# ("S3" and "PYR") and renumbers the ones that follow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Resize(2).EntireRow.Delete() | Out-Null

# Reflect the author's final cell selection recorded in the saved file.
$ws.Range("M9").Select() | Out-Null
